$wb = $excel.ActiveWorkbook

# --- Hoja1: update the daily conversion text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 4.38 = 17092.95 pesos`n✅ 17092.95 pesos = 4.37 = 956.16 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- tasas: update the transfi rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 228.176
$ws2.Range("O10").Value = 3900.2
$ws2.Range("N12").Value = 3915
$ws2.Range("O12").Value = 219
